# Driver.xlsx update:
#  - B2: "HomePage" -> "NewsArticle"
#  - C2: "No" -> "Yes"
#  - C4: "Yes" -> "No"
#  - Active selection on the Driver sheet moves from B11 to D9
# (The now-unused "HomePage" shared string is dropped automatically and the
#  new "NewsArticle" string is appended, which shifts the other shared
#  string indices - that bookkeeping is handled by Excel itself.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Driver")

$ws.Range("B2").Value = "NewsArticle"
$ws.Range("C2").Value = "Yes"
$ws.Range("C4").Value = "No"

$null = $ws.Range("D9").Select()
